$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2021" column (O) to the right of the existing "2020" column (N) ---

# Header row (row 4): copy format from N4 into O4, then set the year value.
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2021

# Bottom border row (row 3): copy the thin/empty bordered cell format from N3 into O3.
$ws.Range("N3").Copy($ws.Range("O3"))

# Data row (row 5): copy N5's format (border + vertical-center) into the new O5 cell,
# update the existing N5 value, and set the new O5 value.
$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("N5").Value = 3.1
$ws.Range("O5").Value = 4.0999999999999996

# Revise an existing data point in the same row.
$ws.Range("L5").Value = 1.6

# Move the active selection to reflect where the editor ended up (was P6, now P4).
$ws.Range("P4").Select()
